$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for all data rows (2-143)
for ($r = 2; $r -le 143; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# Row 2 specific updates: Signalarter (I2) and Alla arter (Q2) counts
$ws.Range("I2").Value = 10
$ws.Range("Q2").Value = 24

# Row 2 species list (R2) gains two new species in alphabetical order
$newText = "Knärot`r`nKoppartaggsvamp`r`nRynkskinn`r`nEntita`r`nGranticka`r`nGränsticka`r`nKandelabersvamp`r`nLeptoporus mollis`r`nMindre hackspett`r`nSpillkråka`r`nTallticka`r`nUllticka`r`nBlodticka`r`nBlomkålssvamp`r`nBlåmossa`r`nDropptaggsvamp`r`nGrovticka`r`nRödgul trumpetsvamp`r`nSkarp dropptaggsvamp`r`nSmal svampklubba`r`nSvart trolldruva`r`nTrådticka`r`nFläcknycklar`r`nBlåsippa"
$ws.Range("R2").Value = $newText
